# Daily price-data refresh: insert a new "today" row at the top of the
# table (row 2, just below the header) with the latest date and the same
# commodity prices, pushing all earlier rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (first data row);
# existing rows 2..30 shift down to 3..31.
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits formatting from the header row above
# it (bold/centered/bordered). Reset it back to the plain "Normal" style
# used by every other data row before writing values.
$ws.Range("A2:D2").Style = "Normal"

# Write the date as literal text (not an auto-converted date serial) by
# briefly using a text number format, matching how the other date cells
# in column A are stored as plain strings.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-20"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Drop the number-format override so the cell ends up with the default
# (unstyled) formatting like its siblings.
$ws.Range("A2").Style = "Normal"
